$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the content: replace "Berufsgenossenschaft Energie Textil Elektro
# Medienerzeugnisse" / "BG ETEM" with "Berufsgenossenschaft Nahrungsmittel" / "BGN"
$ws.Range("B2").Value = "Berufsgenossenschaft Nahrungsmittel"
$ws.Range("B3").Value = "BGN"

# Update the selection to match the new active range
$ws.Range("B2:B3").Select()
